$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value2 = 130734156
$ws.Range("B9").Value2 = 91828
$ws.Range("E9").Value2 = 5432
$ws.Range("F9").Value2 = "Granticka"
$ws.Range("G9").Value2 = "Porodaedalea chrysoloma s.lat."
$ws.Range("H9").Value2 = $null
$ws.Range("K9").Value2 = "teleomorf"
$ws.Range("M9").Value2 = $null
$ws.Range("Q9").Value2 = 443886
$ws.Range("R9").Value2 = 7053279
$ws.Range("AC9").Value2 = "Gott om fruktkroppar i en granhögstubbe."
$ws.Range("AI9").Value2 = $null
$ws.Range("AM9").Value2 = "Stående död trädstam/högstubbe"
$ws.Range("AO9").Value2 = "Standing dead tree/snags # Picea abies"

# Row 10
$ws.Range("A10").Value2 = 130741302
$ws.Range("B10").Value2 = 57884
$ws.Range("E10").Value2 = 100109
$ws.Range("F10").Value2 = "Tretåig hackspett"
$ws.Range("G10").Value2 = "Picoides tridactylus"
$ws.Range("H10").Value2 = "(Linnaeus, 1758)"
$ws.Range("K10").Value2 = $null
$ws.Range("P10").Value2 = "nästsjön, Jmt"
$ws.Range("Q10").Value2 = 443997
$ws.Range("R10").Value2 = 7053008
$ws.Range("AC10").Value2 = "Ringhack färska och äldre"
$ws.Range("AH10").Value2 = $null
$ws.Range("AJ10").Value2 = $null
$ws.Range("AK10").Value2 = $null
$ws.Range("AM10").Value2 = $null
$ws.Range("AO10").Value2 = $null
$ws.Range("AW10").Value2 = "Benny Öwre"
$ws.Range("AX10").Value2 = "Benny Öwre"

# Row 11
$ws.Range("A11").Value2 = 130741291
$ws.Range("Q11").Value2 = 443794
$ws.Range("R11").Value2 = 7053002
$ws.Range("AC11").Value2 = "Ringhack äldre"

# Row 12
$ws.Range("A12").Value2 = 130734145
$ws.Range("M12").Value2 = "färska spår"
$ws.Range("P12").Value2 = "Henrikmyren, Jmt"
$ws.Range("Q12").Value2 = 443967
$ws.Range("R12").Value2 = 7053469
$ws.Range("AC12").Value2 = "Ringhack, färska och äldre, på stambasen av en gran."
$ws.Range("AH12").Value2 = "Granskog"
$ws.Range("AI12").Value2 = "Äldre flerskiktad grandominerad skog."
$ws.Range("AJ12").Value2 = "gran"
$ws.Range("AK12").Value2 = "Picea abies"
$ws.Range("AM12").Value2 = "Trädstam på levande träd"
$ws.Range("AO12").Value2 = "Stem on living tree # Picea abies"
$ws.Range("AW12").Value2 = "Kristian Zackrisson"
$ws.Range("AX12").Value2 = "Kristian Zackrisson"

# Row 13
$ws.Range("A13").Value2 = 130741279
$ws.Range("B13").Value2 = 57884
$ws.Range("E13").Value2 = 100109
$ws.Range("F13").Value2 = "Tretåig hackspett"
$ws.Range("G13").Value2 = "Picoides tridactylus"
$ws.Range("H13").Value2 = "(Linnaeus, 1758)"
$ws.Range("Q13").Value2 = 444007
$ws.Range("R13").Value2 = 7053457
$ws.Range("AC13").Value2 = "Ringhack"

# Row 14
$ws.Range("A14").Value2 = 130741320
$ws.Range("B14").Value2 = 89193
$ws.Range("E14").Value2 = 510
$ws.Range("F14").Value2 = "Doftskinn"
$ws.Range("G14").Value2 = "Cystostereum murrayi"
$ws.Range("H14").Value2 = "(Berk. & M.A.Curtis.) Pouzar"
$ws.Range("P14").Value2 = "nästsjön, Jmt"
$ws.Range("Q14").Value2 = 444326
$ws.Range("R14").Value2 = 7053588
$ws.Range("AC14").Value2 = $null
$ws.Range("AH14").Value2 = $null
$ws.Range("AJ14").Value2 = $null
$ws.Range("AK14").Value2 = $null
$ws.Range("AM14").Value2 = $null
$ws.Range("AO14").Value2 = $null
$ws.Range("AW14").Value2 = "Benny Öwre"
$ws.Range("AX14").Value2 = "Benny Öwre"

# Row 15
$ws.Range("A15").Value2 = 130734180
$ws.Range("Q15").Value2 = 443846
$ws.Range("R15").Value2 = 7053252
$ws.Range("AC15").Value2 = "Växer på grenar av en stående död gran med full längd (bhd ca 30 cm)."
$ws.Range("AM15").Value2 = "Stående död trädstam/högstubbe"
$ws.Range("AO15").Value2 = "Standing dead tree/snags # Picea abies"

# Row 16
$ws.Range("A16").Value2 = 130734170
$ws.Range("Q16").Value2 = 444133
$ws.Range("R16").Value2 = 7053693
$ws.Range("AC16").Value2 = $null

# Row 17
$ws.Range("A17").Value2 = 130734164
$ws.Range("B17").Value2 = 79243
$ws.Range("E17").Value2 = 6425
$ws.Range("F17").Value2 = "Garnlav"
$ws.Range("G17").Value2 = "Alectoria sarmentosa"
$ws.Range("H17").Value2 = "(Ach.) Ach."
$ws.Range("P17").Value2 = "Henrikmyren, Jmt"
$ws.Range("Q17").Value2 = 444362
$ws.Range("R17").Value2 = 7054079
$ws.Range("AC17").Value2 = "På gamla granar."
$ws.Range("AH17").Value2 = "Granskog"
$ws.Range("AJ17").Value2 = "gran"
$ws.Range("AK17").Value2 = "Picea abies"
$ws.Range("AM17").Value2 = "Gren på levande träd"
$ws.Range("AO17").Value2 = "Branch on living tree # Picea abies"
$ws.Range("AW17").Value2 = "Kristian Zackrisson"
$ws.Range("AX17").Value2 = "Kristian Zackrisson"

# Row 35
$ws.Range("A35").Value2 = 130741289
$ws.Range("B35").Value2 = 57884
$ws.Range("E35").Value2 = 100109
$ws.Range("F35").Value2 = "Tretåig hackspett"
$ws.Range("G35").Value2 = "Picoides tridactylus"
$ws.Range("H35").Value2 = "(Linnaeus, 1758)"
$ws.Range("P35").Value2 = "nästsjön, Jmt"
$ws.Range("Q35").Value2 = 443999
$ws.Range("R35").Value2 = 7053066
$ws.Range("AC35").Value2 = "Ringhack äldre"
$ws.Range("AH35").Value2 = $null
$ws.Range("AJ35").Value2 = $null
$ws.Range("AK35").Value2 = $null
$ws.Range("AM35").Value2 = $null
$ws.Range("AO35").Value2 = $null
$ws.Range("AW35").Value2 = "Benny Öwre"
$ws.Range("AX35").Value2 = "Benny Öwre"

# Row 36
$ws.Range("A36").Value2 = 130741294
$ws.Range("B36").Value2 = 57884
$ws.Range("E36").Value2 = 100109
$ws.Range("F36").Value2 = "Tretåig hackspett"
$ws.Range("G36").Value2 = "Picoides tridactylus"
$ws.Range("H36").Value2 = "(Linnaeus, 1758)"
$ws.Range("P36").Value2 = "nästsjön, Jmt"
$ws.Range("Q36").Value2 = 443835
$ws.Range("R36").Value2 = 7052933
$ws.Range("AC36").Value2 = "Ringhack färska och äldre"
$ws.Range("AH36").Value2 = $null
$ws.Range("AI36").Value2 = $null
$ws.Range("AJ36").Value2 = $null
$ws.Range("AK36").Value2 = $null
$ws.Range("AM36").Value2 = $null
$ws.Range("AO36").Value2 = $null
$ws.Range("AW36").Value2 = "Benny Öwre"
$ws.Range("AX36").Value2 = "Benny Öwre"

# Row 37
$ws.Range("A37").Value2 = 130734174
$ws.Range("Q37").Value2 = 444078
$ws.Range("R37").Value2 = 7053449
$ws.Range("AC37").Value2 = "Goda mängder garnlav på flera granar."

# Row 38
$ws.Range("A38").Value2 = 130734189
$ws.Range("B38").Value2 = 79243
$ws.Range("E38").Value2 = 6425
$ws.Range("F38").Value2 = "Garnlav"
$ws.Range("G38").Value2 = "Alectoria sarmentosa"
$ws.Range("H38").Value2 = "(Ach.) Ach."
$ws.Range("P38").Value2 = "Henrikmyren, Jmt"
$ws.Range("Q38").Value2 = 444342
$ws.Range("R38").Value2 = 7053646
$ws.Range("AC38").Value2 = "Relativt rikligt med garnlavsbålar på gran."
$ws.Range("AH38").Value2 = "Granskog"
$ws.Range("AI38").Value2 = "Äldre flerskiktad grandominerad skog med björk och inslag av sälg."
$ws.Range("AJ38").Value2 = "gran"
$ws.Range("AK38").Value2 = "Picea abies"
$ws.Range("AM38").Value2 = "Gren på levande träd"
$ws.Range("AO38").Value2 = "Branch on living tree # Picea abies"
$ws.Range("AW38").Value2 = "Kristian Zackrisson"
$ws.Range("AX38").Value2 = "Kristian Zackrisson"

# Row 39
$ws.Range("A39").Value2 = 130734176
$ws.Range("B39").Value2 = 79243
$ws.Range("E39").Value2 = 6425
$ws.Range("F39").Value2 = "Garnlav"
$ws.Range("G39").Value2 = "Alectoria sarmentosa"
$ws.Range("H39").Value2 = "(Ach.) Ach."
$ws.Range("P39").Value2 = "Henrikmyren, Jmt"
$ws.Range("Q39").Value2 = 444149
$ws.Range("R39").Value2 = 7053319
$ws.Range("AC39").Value2 = "På gran i gles äldre granskog."
$ws.Range("AH39").Value2 = "Granskog"
$ws.Range("AJ39").Value2 = "gran"
$ws.Range("AK39").Value2 = "Picea abies"
$ws.Range("AM39").Value2 = "Gren på levande träd"
$ws.Range("AO39").Value2 = "Branch on living tree # Picea abies"
$ws.Range("AW39").Value2 = "Kristian Zackrisson"
$ws.Range("AX39").Value2 = "Kristian Zackrisson"

# Row 55
$ws.Range("A55").Value2 = 130741311
$ws.Range("B55").Value2 = 91828
$ws.Range("E55").Value2 = 5432
$ws.Range("F55").Value2 = "Granticka"
$ws.Range("G55").Value2 = "Porodaedalea chrysoloma s.lat."
$ws.Range("H55").Value2 = $null
$ws.Range("Q55").Value2 = 443991
$ws.Range("R55").Value2 = 7052955
$ws.Range("AC55").Value2 = $null

# Row 56
$ws.Range("A56").Value2 = 130741293
$ws.Range("B56").Value2 = 57884
$ws.Range("E56").Value2 = 100109
$ws.Range("F56").Value2 = "Tretåig hackspett"
$ws.Range("G56").Value2 = "Picoides tridactylus"
$ws.Range("H56").Value2 = "(Linnaeus, 1758)"
$ws.Range("Q56").Value2 = 443827
$ws.Range("R56").Value2 = 7052932
$ws.Range("AC56").Value2 = "Ringhack äldre"

# Row 62
$ws.Range("A62").Value2 = 130734155
$ws.Range("B62").Value2 = 91804
$ws.Range("E62").Value2 = 1108
$ws.Range("F62").Value2 = "Harticka"
$ws.Range("G62").Value2 = "Pelloporus leporinus"
$ws.Range("H62").Value2 = "(Fr.) Krieglst."
$ws.Range("K62").Value2 = "teleomorf"
$ws.Range("Q62").Value2 = 443861
$ws.Range("R62").Value2 = 7053306
$ws.Range("AC62").Value2 = "Flera fruktkroppar i stambasen av en levande relativt grov gran."
$ws.Range("AM62").Value2 = "Trädstam på levande träd"
$ws.Range("AO62").Value2 = "Stem on living tree # Picea abies"

# Row 63
$ws.Range("A63").Value2 = 130741308
$ws.Range("B63").Value2 = 91828
$ws.Range("E63").Value2 = 5432
$ws.Range("F63").Value2 = "Granticka"
$ws.Range("G63").Value2 = "Porodaedalea chrysoloma s.lat."
$ws.Range("H63").Value2 = $null
$ws.Range("P63").Value2 = "nästsjön, Jmt"
$ws.Range("Q63").Value2 = 444020
$ws.Range("R63").Value2 = 7053458
$ws.Range("AH63").Value2 = $null
$ws.Range("AJ63").Value2 = $null
$ws.Range("AK63").Value2 = $null
$ws.Range("AM63").Value2 = $null
$ws.Range("AO63").Value2 = $null
$ws.Range("AW63").Value2 = "Benny Öwre"
$ws.Range("AX63").Value2 = "Benny Öwre"

# Row 64
$ws.Range("A64").Value2 = 130741283
$ws.Range("B64").Value2 = 57884
$ws.Range("E64").Value2 = 100109
$ws.Range("F64").Value2 = "Tretåig hackspett"
$ws.Range("G64").Value2 = "Picoides tridactylus"
$ws.Range("H64").Value2 = "(Linnaeus, 1758)"
$ws.Range("Q64").Value2 = 443927
$ws.Range("R64").Value2 = 7052967
$ws.Range("AC64").Value2 = "Ringhack äldre"

# Row 65
$ws.Range("A65").Value2 = 130734172
$ws.Range("B65").Value2 = 79243
$ws.Range("E65").Value2 = 6425
$ws.Range("F65").Value2 = "Garnlav"
$ws.Range("G65").Value2 = "Alectoria sarmentosa"
$ws.Range("H65").Value2 = "(Ach.) Ach."
$ws.Range("K65").Value2 = $null
$ws.Range("Q65").Value2 = 444002
$ws.Range("R65").Value2 = 7053503
$ws.Range("AC65").Value2 = $null
$ws.Range("AM65").Value2 = "Gren på levande träd"
$ws.Range("AO65").Value2 = "Branch on living tree # Picea abies"

# Row 66
$ws.Range("A66").Value2 = 130734179
$ws.Range("B66").Value2 = 79243
$ws.Range("E66").Value2 = 6425
$ws.Range("F66").Value2 = "Garnlav"
$ws.Range("G66").Value2 = "Alectoria sarmentosa"
$ws.Range("H66").Value2 = "(Ach.) Ach."
$ws.Range("P66").Value2 = "Henrikmyren, Jmt"
$ws.Range("Q66").Value2 = 443901
$ws.Range("R66").Value2 = 7053212
$ws.Range("AC66").Value2 = $null
$ws.Range("AH66").Value2 = "Granskog"
$ws.Range("AJ66").Value2 = "gran"
$ws.Range("AK66").Value2 = "Picea abies"
$ws.Range("AM66").Value2 = "Gren på levande träd"
$ws.Range("AO66").Value2 = "Branch on living tree # Picea abies"
$ws.Range("AW66").Value2 = "Kristian Zackrisson"
$ws.Range("AX66").Value2 = "Kristian Zackrisson"
